$d = $word.ActiveDocument

$replacements = @(
    @("2024-04-15 Monday", "2024-04-16 Tuesday"),
    @("66×70=4620", "69×28=1932"),
    @("82×52=4264", "17×18=306"),
    @("32×49=1568", "22×45=990"),
    @("26×60=1560", "79×44=3476"),
    @("54×25=1350", "21×63=1323"),
    @("44×58=2552", "21×64=1344"),
    @("16×87=1392", "98×30=2940"),
    @("81×38=3078", "23×89=2047"),
    @("91×86=7826", "69×38=2622"),
    @("96×19=1824", "50×37=1850"),
    @("38×87=3306", "44×36=1584"),
    @("15×92=1380", "32×35=1120"),
    @("77×64=4928", "72×96=6912"),
    @("13×63=819", "92×87=8004"),
    @("41×59=2419", "81×26=2106"),
    @("63×40=2520", "49×58=2842"),
    @("20×40=800", "19×67=1273"),
    @("80×78=6240", "27×79=2133"),
    @("22×89=1958", "34×99=3366"),
    @("17×13=221", "78×40=3120"),
    @("53×38=2014", "90×28=2520"),
    @("69×37=2553", "34×40=1360"),
    @("78×45=3510", "69×67=4623"),
    @("12×40=480", "33×30=990"),
    @("40×39=1560", "49×65=3185")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
